# Reahl_2020 ALLDATA.xlsx — "Change climate color on ancient sample,
# add climate supp. figures"
#
# The only data-level edit is on the "Freq. of Occurrence" sheet: the
# "climatecolor" swatch for the ancient-sample row 115 changes from the
# yellow hex code (#F0E442) to black (#000000). The sheet's last-saved
# view/selection state also moved (scrolled further down, new active
# cell) as a side effect of the author's editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Freq. of Occurrence")
$ws.Activate()

# F115: climatecolor #F0E442 -> #000000
$ws.Range("F115").Value = "#000000"

# Sheet view moved: topLeftCell A50 -> A96, selection J69 -> I118
$excel.ActiveWindow.ScrollRow = 96
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I118").Select()
